$wb = $excel.ActiveWorkbook

# Sheet 1: "HSV Log OLS"
$ws1 = $wb.Worksheets.Item("HSV Log OLS")
$ws1.Range("A2").Value = 0.127
$ws1.Range("B2").Value = 3.893
$ws1.Range("C2").Value = 0.874

# Sheet 2: "HSV PPML"
$ws2 = $wb.Worksheets.Item("HSV PPML")
$ws2.Range("A2").Value = 0.023
$ws2.Range("B2").Value = 1.145
$ws2.Range("C2").Value = 0.874

# Sheet 3: "HSVT NLLSQ"
$ws3 = $wb.Worksheets.Item("HSVT NLLSQ")
$ws3.Range("A2").Value = -0.051
$ws3.Range("B2").Value = 0.402
$ws3.Range("C2").Value = 15214.24
$ws3.Range("D2").Value = 0.119
